# Edit script for "STOCK MARKET INVESTMENT STRATEGY.pptx"
# Adds a commentary textbox under the chart on six slides and repositions
# the chart picture on each of those slides to make room for the text.

$p = $ppt.ActivePresentation

function Set-ChartPicture {
    param($shape, [double]$left, [double]$top, [double]$width, [double]$height)
    $shape.Left = $left
    $shape.Top = $top
    $shape.Width = $width
    $shape.Height = $height
}

function New-CommentBox {
    param($slide, [double]$left, [double]$top, [double]$width, [double]$height, [string]$name)
    $tb = $slide.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $tb.Name = $name
    $tb.Fill.Visible = $false
    $tb.Line.Visible = $false
    $tb.TextFrame.WordWrap = -1
    $tb.TextFrame.AutoSize = 0
    $tb.TextFrame.MarginLeft = 0
    $tb.TextFrame.MarginRight = 0
    $tb.TextFrame.MarginTop = 0
    $tb.TextFrame.MarginBottom = 0
    $tb.TextFrame.VerticalAnchor = 1
    return $tb
}

# ---------------------------------------------------------------------------
# Slide 10 - "30-day Moving Average"
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$pic10 = $s10.Shapes.Item(3)
Set-ChartPicture $pic10 52.363622 111.009449 850.181811 317.981102

$tb10 = New-CommentBox $s10 52.363543 435.118110 850.181811 104.881890 "Google Shape;286;p27"
$tb10.TextFrame.TextRange.Text = " - moving average for Microsoft is showing an uptrend so it is most likely to maintain that trend given its investments and its products which shows dominancy in the computer space and its expansion in the technological space."
$tb10.TextFrame.TextRange.Font.Size = 14

# ---------------------------------------------------------------------------
# Slide 11 - "Bitcoin: 10-day Moving Average"
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)

$tb11 = New-CommentBox $s11 62.647717 419.614173 808.852205 104.881890 "Google Shape;293;p28"
$tb11.TextFrame.TextRange.Text = " - rolling sum of 10 days due to not having enough data and using a rolling sum of higher days can be misleading but we can see that the moving average is going down since 2021 as we mentioned that bitcoin has lost value in the past 2 years."
$tb11.TextFrame.TextRange.Font.Size = 14

# ---------------------------------------------------------------------------
# Slide 12 - "Amazon Risk vs Expected Return"
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$pic12 = $s12.Shapes.Item(3)
Set-ChartPicture $pic12 271.272677 112.586142 437.333701 291.050394

$tb12 = New-CommentBox $s12 74.284094 427.413858 873.354331 104.881890 "Google Shape;300;p29"
$tb12.TextFrame.TextRange.Text = " -The expected return to Risk is reasonable considering that we are looking to invest long-term."
$tb12.TextFrame.TextRange.Font.Size = 14

# ---------------------------------------------------------------------------
# Slide 13 - "Microsoft Risk vs Expected Return"
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$pic13 = $s13.Shapes.Item(3)
Set-ChartPicture $pic13 294.545433 119.917323 429.375669 282.264488

$tb13 = New-CommentBox $s13 57.556772 435.118110 873.354331 104.881890 "Google Shape;307;p30"
$tb13.TextFrame.TextRange.Text = " - The expected return to Risk is reasonable considering that we are looking to invest long term."
$tb13.TextFrame.TextRange.Font.Size = 14

# ---------------------------------------------------------------------------
# Slide 14 - "Bitcoin Risk vs Expected Return"
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$pic14 = $s14.Shapes.Item(3)
Set-ChartPicture $pic14 253.818110 105.625984 455.807795 291.464961

$tb14 = New-CommentBox $s14 90.181732 420.319685 800.727323 104.881890 "Google Shape;314;p31"
$fullText14 = "- Bitcoin has the highest risk as to the other stocks and the returns are a loss so it would not be wise to invest to it due to the votility of the stock and it not yielding enough returns for a profit."
$tb14.TextFrame.TextRange.Text = $fullText14
$tb14.TextFrame.TextRange.Font.Size = 14
$volStart = $fullText14.IndexOf("votility") + 1
$volLen = "votility".Length
$volRange = $tb14.TextFrame.TextRange.Characters($volStart, $volLen)
$volRange.Font.Size = 14

# ---------------------------------------------------------------------------
# Slide 15 - "Inflation Rate"
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$pic15 = $s15.Shapes.Item(3)
Set-ChartPicture $pic15 98.181890 111.009449 770.909134 308.626929

$tb15 = New-CommentBox $s15 88.727165 428.990551 780.363858 104.881890 "Google Shape;321;p32"
$tb15.TextFrame.TextRange.Text = " - The overall trend is increasing even though there was a drastic drop in 2008. In 2020 we observed a sharp increase in inflation and it has never shown any sign of dropping."
$tb15.TextFrame.TextRange.Font.Size = 14
